$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.861.18'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.01%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.516.87'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.11%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '601.22'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.89%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '180.93'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.47%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.516.68'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.13%  '
$ws.Range("E9").Value = '  +0.13%  '
$ws.Range("E10").Value = '  +6.64%  '
$ws.Range("E11").Value = '  -1.64%  '
$ws.Range("E12").Value = '  +0.41%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.128.19'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.15%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.23'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +10.05%  '
$ws.Range("E15").Value = '  +0.96%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '67.876.23'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.09%  '
$ws.Range("E17").Value = '  +0.40%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.514.39'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.12%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.37'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.32%  '
$ws.Range("E20").Value = '  +2.08%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '400.45'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.87%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '8.00'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.65%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '73.79'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.98%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.544'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.07%  '
$ws.Range("E25").Value = '  +0.35%  '
$ws.Range("E26").Value = '  +0.82%  '
$ws.Range("E27").Value = '  +0.70%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.49'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.11%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.178'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.50%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.997'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.03%  '
$ws.Range("E31").Value = '  -0.52%  '
$ws.Range("E32").Value = '  -0.86%  '
$ws.Range("E33").Value = '  +1.06%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '23.93'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.16%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '7.51'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.56%  '
$ws.Range("E36").Value = '  +0.11%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.64'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.14%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '163.31'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.35%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.881'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.72%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.92'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.48%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.81'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +8.74%  '
$ws.Range("E42").Value = '  -1.73%  '
$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.904.15'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.30%  '
$ws.Range("B44").Value = 'Filecoin'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.71'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.14%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '26.42'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.66%  '
$ws.Range("E46").Value = '  -2.21%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '26.95'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.97%  '
$ws.Range("E48").Value = '  -1.08%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '351.58'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.48%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0305'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.20%  '
$ws.Range("E51").Value = '  -1.13%  '
